# Updated symbol list on Sat Dec 17 20:33:05 UTC 2022 with GitHub Actions
#
# The "Price" (D) and "Volume(1h)" (E) columns are refreshed with new
# scraped values, and two adjacent rows (BKEXToken / CEJI) swap rank so
# their Coin name / Link / Price / Volume cells trade places.
#
# Every one of these cells holds plain TEXT in the workbook (not a real
# number), even the price figures. Writing a numeric-looking string via
# .Value would make Excel re-interpret it as a number (dropping
# significant trailing/leading zeros, e.g. "0.05570" -> 0.0557), so each
# numeric-looking price is entered with a leading apostrophe to force a
# text entry - exactly what a person retyping these figures in Excel
# would need to do. Excel then marks that cell "quote prefixed"
# (number format style change); resetting the cell style back to
# "Normal" immediately afterwards keeps the cell's look/format exactly
# as it was before, with only the text content changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextPrice {
    param($addr, $value)
    $ws.Range($addr).Value = "'" + $value
    $ws.Range($addr).Style = "Normal"
}

# Row 2 - BNB
Set-TextPrice "D2" "235.34"

# Row 3 - OKB
Set-TextPrice "D3" "21.68"

# Row 4 - HuobiToken
Set-TextPrice "D4" "5.364"

# Row 5 - Cronos
Set-TextPrice "D5" "0.05570"

# Row 6 - GateToken
Set-TextPrice "D6" "3.367"

# Row 7 - KuCoinToken
Set-TextPrice "D7" "6.461"

# Row 8 - MXToken
Set-TextPrice "D8" "0.8011"

# Row 9 - FTXToken
Set-TextPrice "D9" "1.039"

# Row 10 - WazirX
Set-TextPrice "D10" "0.1405"

# Row 11 - MandalaExchangeToken
Set-TextPrice "D11" "0.07249"

# Row 12 - LiechtensteinCryptoassetsExchange
Set-TextPrice "D12" "0.03183"

# Row 13 - BitrueCoin
Set-TextPrice "D13" "0.02937"

# Row 14 - BitMartToken
Set-TextPrice "D14" "0.09247"

# Row 15 - BitForexToken
Set-TextPrice "D15" "0.001660"

# Row 16 - MCDex
Set-TextPrice "D16" "3.257"

# Row 17 - CoinExToken
Set-TextPrice "D17" "0.04771"

# Row 18 - One
Set-TextPrice "D18" "0.0005713"
$ws.Range("E18").Value = "17OneONEWorstin24h"

# Row 19 - TigerCash
Set-TextPrice "D19" "0.006261"

# Row 20
Set-TextPrice "D20" "0.005071"

# Row 22
Set-TextPrice "D22" "0.0001501"

# Row 23
Set-TextPrice "D23" "0.0004202"

# Row 24 - LEO
Set-TextPrice "D24" "3.938"
$ws.Range("E24").Value = "23LEOLEOBestin24h"

# Row 25
Set-TextPrice "D25" "2.201"

# Row 27
Set-TextPrice "D27" "0.1307"

# Row 40 - IDEX
Set-TextPrice "D40" "0.04117"

# Row 41 - KickToken
Set-TextPrice "D41" "0.007036"

# Rows 42 / 43 swapped rank: BKEXToken now ranks above CEJI
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextPrice "D42" "0.1039"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextPrice "D43" "0.002922"
$ws.Range("E43").Value = "42CEJICEJI"

# Row 44 - LocalTraders
Set-TextPrice "D44" "0.008884"

# Row 45 - CoinLion
Set-TextPrice "D45" "0.00005439"

# Row 47 - CoinbaseStockToken
Set-TextPrice "D47" "0.6803"

# Row 48 - BOLO
Set-TextPrice "D48" "0.03274"
$ws.Range("E48").Value = "47BOLOBOLO"

# Row 49 - CryptobidCoin
Set-TextPrice "D49" "0.00002101"
